$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns: F (计划开始时间 / planned start time) and G (计划完成时间 / planned finish time)
$ws.Range("F1").Value = "计划开始时间"
$ws.Range("G1").Value = "计划完成时间"

# Planned start / finish date serials (as captured from the source data), rows 2-16
$startSerial = @(42399, 42399, 42399, 42459, 42459, 42459, 42459, 42459, 42459, 42459, 42459, 42459, 42459, 42459, 42459)
$finishSerial = @(42460, 42460, 42460, 42734, 42734, 42734, 42734, 42734, 42734, 42734, 42734, 42734, 42734, 42734, 42734)

for ($i = 0; $i -lt 15; $i++) {
    $row = $i + 2
    $fCell = $ws.Cells.Item($row, 6)
    $gCell = $ws.Cells.Item($row, 7)
    $fCell.Value = $startSerial[$i]
    $gCell.Value = $finishSerial[$i]
    $fCell.NumberFormat = "mm-dd-yy"
    $gCell.NumberFormat = "mm-dd-yy"
}

# Column widths for the new date columns
$ws.Columns.Item(6).ColumnWidth = 14.857142857142858
$ws.Columns.Item(7).ColumnWidth = 11.714285714285714

# --- Column B ("状态"/status) switches from a percentage display to a plain number
$ws.Range("B1:B16").NumberFormat = "General"
$ws.Columns.Item(2).ColumnWidth = 8.285714285714286

# --- Selection moves to G18, matching the new filterable/date-aware layout
$ws.Range("G18").Select()
